$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 20.76103376777028
$ws.Range("D2").Value = 100.5891923750444
